$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New stimuli order: for each data row (2..33) set the numeric seed/index (col B),
# the image path (col C), the word (col D) and the category (col E).
# Column A (row index 0..31) is unchanged.

$rows = @(
    @{B=76;  C="face/face095.png"; D="drehen";    E="face"},
    @{B=105; C="car/car121.png";  D="husten";    E="car"},
    @{B=82;  C="face/face097.png"; D="enden";     E="face"},
    @{B=72;  C="face/face096.png"; D="fühlen";    E="face"},
    @{B=48;  C="car/car096.png";  D="pflegen";   E="car"},
    @{B=88;  C="face/face110.png"; D="füttern";   E="face"},
    @{B=116; C="face/face117.png"; D="regnen";    E="face"},
    @{B=86;  C="car/car089.png";  D="bitten";    E="car"},
    @{B=1;   C="face/face078.png"; D="rücken";    E="face"},
    @{B=89;  C="car/car117.png";  D="tagen";     E="car"},
    @{B=6;   C="car/car094.png";  D="schicken";  E="car"},
    @{B=66;  C="car/car110.png";  D="wiegen";    E="car"},
    @{B=61;  C="face/face086.png"; D="drohen";    E="face"},
    @{B=12;  C="car/car070.png";  D="backen";    E="car"},
    @{B=101; C="face/face091.png"; D="hoffen";    E="face"},
    @{B=47;  C="face/face080.png"; D="nehmen";    E="face"},
    @{B=78;  C="face/face098.png"; D="dauern";    E="face"},
    @{B=69;  C="car/car085.png";  D="tauschen";  E="car"},
    @{B=121; C="face/face108.png"; D="runden";    E="face"},
    @{B=8;   C="face/face077.png"; D="rasen";     E="face"},
    @{B=92;  C="face/face068.png"; D="fesseln";   E="face"},
    @{B=2;   C="car/car077.png";  D="ehren";     E="car"},
    @{B=28;  C="car/car088.png";  D="raten";     E="car"},
    @{B=108; C="face/face112.png"; D="scheitern"; E="face"},
    @{B=59;  C="car/car107.png";  D="schenken";  E="car"},
    @{B=4;   C="car/car092.png";  D="klappen";   E="car"},
    @{B=44;  C="car/car069.png";  D="biegen";    E="car"},
    @{B=43;  C="car/car101.png";  D="posten";    E="car"},
    @{B=115; C="car/car091.png";  D="antun";     E="car"},
    @{B=87;  C="car/car095.png";  D="bleiben";   E="car"},
    @{B=113; C="face/face074.png"; D="stechen";   E="face"},
    @{B=18;  C="face/face076.png"; D="haken";     E="face"}
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $r++
}
